$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.702.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.704.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9976"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9976"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.515"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.18"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08814"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.358"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +11.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001326"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.552"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.702.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07126"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.775"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9968"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.688.68"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.002"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.309"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.120"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.79"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.417"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +26.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.887.35"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.15%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08699"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.383"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +18.84%  "

$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.947"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.37%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.97"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2734"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02783"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09008"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.482"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7694"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7219"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.173"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9969"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.90"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.315"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000378"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.22%  "
